$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates reflecting refreshed crypto price/volume data.
# Numeric-looking text values (e.g. "517.93") are protected with a temporary
# Text number format so Excel does not silently convert them to real numbers;
# the style is reset back to "Normal" immediately after so no stray number
# format is left applied to the cell.

$ws.Range("D2").Value = "58.051.00"
$ws.Range("E2").Value = "  -1.71%  "
$ws.Range("D3").Value = "2.465.66"
$ws.Range("E3").Value = "  -1.86%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "517.93"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.93%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "133.27"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.01%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.557"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.88%  "
$ws.Range("D9").Value = "2.475.71"
$ws.Range("E9").Value = "  -1.67%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0978"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.78%  "
$ws.Range("E11").Value = "  -0.80%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.28"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.52%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.336"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.90%  "
$ws.Range("D14").Value = "2.906.41"
$ws.Range("E14").Value = "  -1.83%  "
$ws.Range("D15").Value = "57.991.02"
$ws.Range("E15").Value = "  -1.61%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.87"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.25%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000134"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.82%  "
$ws.Range("D18").Value = "2.476.31"
$ws.Range("E18").Value = "  -0.97%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.57"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.44%  "
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "319.30"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.53%  "
$ws.Range("B21").Value = "Polkadot"
$ws.Range("C21").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.15"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.42%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.13%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.71"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.88%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.50"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.18%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.408"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.73%  "
$ws.Range("E26").Value = "  -0.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.161"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.87%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.32"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.12%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "170.96"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.11%  "
$ws.Range("D30").Value = "0.0₃0742"
$ws.Range("E30").Value = "  -3.19%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.30"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.94%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.68"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.64%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.17"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.75%  "
$ws.Range("E34").Value = "  +0.00%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.997"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.06%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.01"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.14%  "
$ws.Range("E37").Value = "  -4.04%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.96"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.59%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.56"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.05%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.46"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.06%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.790"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.04%  "
$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.42"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.65%  "
$ws.Range("B43").Value = "Bittensor"
$ws.Range("C43").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "271.86"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.68%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.02"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.05%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.592"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.08%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "123.25"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.48%  "
$ws.Range("E47").Value = "  -1.88%  "
$ws.Range("E48").Value = "  -2.71%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0212"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.79%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "16.90"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.57%  "
$ws.Range("D51").Value = "1.730.48"
$ws.Range("E51").Value = "  -1.66%  "
